# Auto-generated edit script applying numeric updates from the commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 125.333336
$ws.Range("I11").Value = 125.333336
$ws.Range("K11").Value = 125.333336
$ws.Range("M11").Value = 14.666664
# Row 40
$ws.Range("H40").Value = 7228.875
$ws.Range("I40").Value = 6110.5
$ws.Range("K40").Value = 6110.5
$ws.Range("M40").Value = -5935.5
# Row 70
$ws.Range("H70").Value = 8289.6
$ws.Range("I70").Value = 2450
$ws.Range("K70").Value = 7350
$ws.Range("M70").Value = -7080
# Row 73
$ws.Range("H73").Value = 8289.6
$ws.Range("I73").Value = 2450
$ws.Range("K73").Value = 7350
$ws.Range("M73").Value = -6414
# Row 74
$ws.Range("H74").Value = 6698.2856
$ws.Range("I74").Value = 6698.2856
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6698.2856
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -5762.2856
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 6698.2856
$ws.Range("I77").Value = 6698.2856
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 33491.428
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -28811.428
$ws.Range("N77").ClearContents()
# Row 94
$ws.Range("H94").Value = 3701.182
$ws.Range("I94").Value = 3701.182
$ws.Range("K94").Value = 3701.182
$ws.Range("M94").Value = -3250.182
# Row 132
$ws.Range("H132").Value = 6464.3076
$ws.Range("I132").Value = 1526.3334
$ws.Range("J132").Value = 17574.75
$ws.Range("K132").Value = 4579.0002
$ws.Range("L132").Value = 52724.25
$ws.Range("M132").Value = -2049.0002
$ws.Range("N132").Value = -57784.25
# Row 137
$ws.Range("H137").Value = 2440.889
$ws.Range("I137").Value = 1193.6
$ws.Range("K137").Value = 3580.8
$ws.Range("M137").Value = -1030.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 12434.117
$ws.Range("J44").Value = 12434.117
$ws.Range("L44").Value = 12434.117
$ws.Range("N44").Value = -13410.117
# Row 132
$ws.Range("H132").Value = 1969.7878
$ws.Range("I132").Value = 1966.7667
$ws.Range("K132").Value = 5900.300099999999
$ws.Range("M132").Value = -3370.300099999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 967.5
$ws.Range("I16").Value = 964.5454999999999
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 964.5454999999999
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -677.5454999999999
$ws.Range("N16").Value = -1574
# Row 22
$ws.Range("H22").Value = 2256.3333
$ws.Range("J22").Value = 3100
$ws.Range("L22").Value = 3100
$ws.Range("N22").Value = -3800
# Row 99
$ws.Range("H99").Value = 1623.1818
$ws.Range("I99").Value = 1623.1818
$ws.Range("K99").Value = 1623.1818
$ws.Range("M99").Value = -125.1818000000001
# Row 107
$ws.Range("H107").Value = 620.3103599999999
$ws.Range("I107").Value = 583.2632
$ws.Range("J107").Value = 690.7
$ws.Range("K107").Value = 583.2632
$ws.Range("L107").Value = 690.7
$ws.Range("M107").Value = 1336.7368
$ws.Range("N107").Value = -4530.7
# Row 113
$ws.Range("H113").Value = 967.5
$ws.Range("I113").Value = 964.5454999999999
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 964.5454999999999
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1205.4545
$ws.Range("N113").Value = -5340
# Row 126
$ws.Range("H126").Value = 1623.1818
$ws.Range("I126").Value = 1623.1818
$ws.Range("K126").Value = 4869.5454
$ws.Range("M126").Value = -2399.5454
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 3136.6538
$ws.Range("I132").Value = 2891.0454
$ws.Range("K132").Value = 8673.136200000001
$ws.Range("M132").Value = -6143.136200000001
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 113
$ws.Range("H113").Value = 1275.2106
$ws.Range("J113").Value = 1451
$ws.Range("L113").Value = 4353
$ws.Range("N113").Value = -8693
# Row 122
$ws.Range("H122").Value = 619.1111
$ws.Range("I122").Value = 619.1111
$ws.Range("K122").Value = 5571.9999
$ws.Range("M122").Value = -3121.9999
# Row 137
$ws.Range("H137").Value = 5000
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -25200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 28999.834
$ws.Range("J95").Value = 28999.834
$ws.Range("L95").Value = 28999.834
$ws.Range("N95").Value = -34491.834
# Row 122
$ws.Range("H122").Value = 1469.8462
$ws.Range("I122").Value = 1469.8462
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4409.5386
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1959.5386
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 3385.4
$ws.Range("I126").Value = 3291.2307
$ws.Range("K126").Value = 9873.6921
$ws.Range("M126").Value = -7403.6921

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1999.6666
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("N22").Value = -3090
# Row 27
$ws.Range("H27").Value = 1999.6666
$ws.Range("J27").Value = 2500
$ws.Range("L27").Value = 2500
$ws.Range("N27").Value = -2714
# Row 46
$ws.Range("H46").Value = 2727.1904
$ws.Range("I46").Value = 1186.1111
$ws.Range("K46").Value = 1186.1111
$ws.Range("M46").Value = -998.1111000000001
# Row 55
$ws.Range("H55").Value = 812.65
$ws.Range("I55").Value = 868.5333000000001
$ws.Range("J55").Value = 645
$ws.Range("K55").Value = 868.5333000000001
$ws.Range("L55").Value = 645
$ws.Range("M55").Value = -695.5333000000001
$ws.Range("N55").Value = -991
# Row 68
$ws.Range("H68").Value = 8714.571
$ws.Range("I68").Value = 6001
$ws.Range("K68").Value = 6001
$ws.Range("M68").Value = -5252
# Row 71
$ws.Range("H71").Value = 8714.571
$ws.Range("I71").Value = 6001
$ws.Range("K71").Value = 30005
$ws.Range("M71").Value = -26261
# Row 132
$ws.Range("H132").Value = 5612.1577
$ws.Range("I132").Value = 4240.385
$ws.Range("J132").Value = 8584.333000000001
$ws.Range("K132").Value = 12721.155
$ws.Range("L132").Value = 25752.999
$ws.Range("M132").Value = -10191.155
$ws.Range("N132").Value = -30812.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 99
$ws.Range("H99").Value = 63000
$ws.Range("I99").Value = 63000
$ws.Range("K99").Value = 63000
$ws.Range("M99").Value = -60005
# Row 132
$ws.Range("H132").Value = 1315.3889
$ws.Range("I132").Value = 1333.9412
$ws.Range("K132").Value = 4001.8236
$ws.Range("M132").Value = -1471.8236

